{"js": "// The cover letter body is a single paragraph / single run containing plain\n// text split by manual line breaks (<w:br/>), with blank lines represented\n// as two consecutive breaks. We locate each affected \"paragraph\" (the text\n// between one pair of breaks and the next) using short, uniquely-matching\n// anchor phrases - well under Word's search-string length limit - and\n// either replace or delete it, preserving the surrounding break structure.\n\nconst body = context.document.body;\n\n// Find the Range spanning from the start of `startAnchor` to the end of\n// `endAnchor` (inclusive). Both anchors must be unique in the document.\n// If startAnchor === endAnchor a single search result is used directly.\nasync function getAnchoredRange(startAnchor, endAnchor) {\n  const startResults = body.search(startAnchor, { matchCase: true, matchWildcards: false });\n  startResults.load(\"items\");\n  await context.sync();\n  if (startResults.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for start anchor '\" + startAnchor + \"', found \" + startResults.items.length\n    );\n  }\n\n  if (startAnchor === endAnchor) {\n    return startResults.items[0];\n  }\n\n  const endResults = body.search(endAnchor, { matchCase: true, matchWildcards: false });\n  endResults.load(\"items\");\n  await context.sync();\n  if (endResults.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for end anchor '\" + endAnchor + \"', found \" + endResults.items.length\n    );\n  }\n\n  return startResults.items[0].expandTo(endResults.items[0]);\n}\n\n// Replace (or delete, when replacement === \"\") the text running from\n// startAnchor through endAnchor (inclusive of both anchors).\nasync function replaceSpan(startAnchor, endAnchor, replacement) {\n  const range = await getAnchoredRange(startAnchor, endAnchor);\n  range.insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 1) Replace the opening paragraph (old intro) with the entire block of\n//    new/rewritten paragraphs, joined with double manual line-breaks\n//    (\\u000b\\u000b), matching the document's existing paragraph style.\n// ---------------------------------------------------------------------\nconst newParagraphs = [\n  \"I am writing to express my interest in the position of Data Scientist at your company. With a strong background in data analysis and machine learning, I believe I would be a valuable asset to your team.\",\n  \"I have a Bachelor's degree in Technology from NIT Durgapur, where I gained a solid foundation in programming, statistics, and mathematics. During my time at NIT Durgapur, I also completed a Full Stack Data Science Bootcamp at iNeuron, further enhancing my skills in data science and machine learning.\",\n  \"In my previous role as a Senior Associate Data Scientist at Affine, I led a team in scraping data from gaming sites and creating concise dashboards with various charts for insights. I also supervised a data science team in generating and clustering a large analytical dataset and provided real-time visibility in predicting retention and ARPDAU using FB-Prophet. Additionally, I analyzed sentiments from web servers using NLP techniques and built a UI interface for automation tasks.\",\n  \"I have also gained experience in system analysis and implementation during my time as an Assistant System Engineer at Tata Consultancy Services. I conducted system analysis, troubleshooting, and automation using SQL and Unix.\",\n  \"In terms of projects, I have worked on various data science projects, including insurance premium prediction, WhatsApp chat analysis, credit card fraud detection, and stock price prediction. These projects have allowed me to apply my skills in regression, classification, and time series forecasting using various machine learning algorithms.\",\n  \"I possess a strong skill set in Python, machine learning, SQL, NLP, time series forecasting, deep learning, and statistics. I am also proficient in using libraries such as NumPy, Pandas, TensorFlow, Keras, NLTK, and Seaborn. Additionally, I have experience with tools like Looker, Streamlit, PowerBI, and web scraping.\",\n  \"I am a dedicated and hardworking individual with excellent problem-solving and teamwork skills. I am also highly organized and have strong time management abilities.\",\n  \"Outside of work, I have served as the captain of the NIT Durgapur Cricket team at the national level in 2019 and have been recognized with the \\\"Best Employee of the Month\\\" award for my exemplary performance.\"\n];\nconst newCombinedText = newParagraphs.join(\"\\u000b\\u000b\");\n\nawait replaceSpan(\n  \"I am writing to express my interest in\",\n  \"uld be a valuable addition to your team.\",\n  newCombinedText\n);\n\n// ---------------------------------------------------------------------\n// 2) The four remaining old paragraphs (ffine/Affine role, iNeuron\n//    internship, B.Tech/bootcamp, soft-skills) are now superseded by the\n//    block inserted above, so delete each of them along with its full\n//    trailing blank-line break pair (\\u000b\\u000b). The break pair that\n//    precedes each deleted paragraph is left in place, so exactly one\n//    blank line remains between the surviving paragraphs.\n// ---------------------------------------------------------------------\nawait replaceSpan(\n  \"In my current role as a Senior Associate\",\n  \"models and dashboards for game launches.\\u000b\\u000b\",\n  \"\"\n);\nawait replaceSpan(\n  \"During my internship at iNeuron, I worke\",\n  \"alysis, documentation, and user support.\\u000b\\u000b\",\n  \"\"\n);\nawait replaceSpan(\n  \"I hold a B.Tech degree from NIT Durgapur\",\n  \"as NumPy, Pandas, TensorFlow, and Keras.\\u000b\\u000b\",\n  \"\"\n);\nawait replaceSpan(\n  \"In addition to my technical skills, I po\",\n  \"er to learn and adapt to new challenges.\\u000b\\u000b\",\n  \"\"\n);\n\nbody.load(\"text\");\nawait context.sync();\nreturn body.text;\n", "ps1": "# The cover letter body is a single paragraph / single run containing plain\n# text split by manual line breaks (Chr(11), i.e. <w:br/> in the OOXML),\n# with blank lines represented as two consecutive breaks. We locate each\n# affected \"paragraph\" with Find/Execute against $d.Content (which reliably\n# mutates Start/End in place), then either overwrite or clear its .Text,\n# preserving the surrounding break structure.\n\n$d = $word.ActiveDocument\n$vbreak = [char]11\n\nfunction Find-ExactRange($doc, [string]$searchText) {\n    $r = $doc.Content\n    $find = $r.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute()\n    if (-not $found) {\n        throw \"Text not found: $searchText\"\n    }\n    return $r\n}\n\n# ---------------------------------------------------------------------\n# 1) Replace the opening paragraph (old intro) with the entire block of\n#    new/rewritten paragraphs, joined with double manual line-breaks,\n#    matching the document's existing paragraph style.\n# ---------------------------------------------------------------------\n$newParagraphs = @(\n    \"I am writing to express my interest in the position of Data Scientist at your company. With a strong background in data analysis and machine learning, I believe I would be a valuable asset to your team.\",\n    \"I have a Bachelor's degree in Technology from NIT Durgapur, where I gained a solid foundation in programming, statistics, and mathematics. During my time at NIT Durgapur, I also completed a Full Stack Data Science Bootcamp at iNeuron, further enhancing my skills in data science and machine learning.\",\n    \"In my previous role as a Senior Associate Data Scientist at Affine, I led a team in scraping data from gaming sites and creating concise dashboards with various charts for insights. I also supervised a data science team in generating and clustering a large analytical dataset and provided real-time visibility in predicting retention and ARPDAU using FB-Prophet. Additionally, I analyzed sentiments from web servers using NLP techniques and built a UI interface for automation tasks.\",\n    \"I have also gained experience in system analysis and implementation during my time as an Assistant System Engineer at Tata Consultancy Services. I conducted system analysis, troubleshooting, and automation using SQL and Unix.\",\n    \"In terms of projects, I have worked on various data science projects, including insurance premium prediction, WhatsApp chat analysis, credit card fraud detection, and stock price prediction. These projects have allowed me to apply my skills in regression, classification, and time series forecasting using various machine learning algorithms.\",\n    \"I possess a strong skill set in Python, machine learning, SQL, NLP, time series forecasting, deep learning, and statistics. I am also proficient in using libraries such as NumPy, Pandas, TensorFlow, Keras, NLTK, and Seaborn. Additionally, I have experience with tools like Looker, Streamlit, PowerBI, and web scraping.\",\n    \"I am a dedicated and hardworking individual with excellent problem-solving and teamwork skills. I am also highly organized and have strong time management abilities.\",\n    \"Outside of work, I have served as the captain of the NIT Durgapur Cricket team at the national level in 2019 and have been recognized with the `\"Best Employee of the Month`\" award for my exemplary performance.\"\n)\n$blankLine = [string]$vbreak + [string]$vbreak\n$newCombinedText = [string]::Join($blankLine, $newParagraphs)\n\n$oldIntro = \"I am writing to express my interest in the Senior Associate Data Scientist position at your company. With a strong background in data science and a proven track record of delivering impactful insights and solutions, I believe I would be a valuable addition to your team.\"\n$r = Find-ExactRange $d $oldIntro\n$r.Text = $newCombinedText\n\n# ---------------------------------------------------------------------\n# 2) The four remaining old paragraphs (ffine/Affine role, iNeuron\n#    internship, B.Tech/bootcamp, soft-skills) are now superseded by the\n#    block inserted above, so delete each of them along with its full\n#    trailing blank-line break pair. The break pair preceding each\n#    deleted paragraph is left in place, so exactly one blank line\n#    remains between the surviving paragraphs.\n# ---------------------------------------------------------------------\n$oldParasToRemove = @(\n    \"In my current role as a Senior Associate Data Scientist at ffine, I have led a team in scraping data from multiple gaming sites and creating concise dashboards with insightful charts using Google Sheets. I have also supervised a data science team in generating and clustering analytical datasets, providing real-time visibility to leadership for decision-making. Additionally, I have analyzed sentiments from web servers and designed live revenue models and dashboards for game launches.\",\n    \"During my internship at iNeuron, I worked on an insurance premium prediction project where I conducted exploratory data analysis, implemented regression models, and deployed the project as a web application on AWS. I also have experience as an Assistant System Engineer at Tata Consultancy Services, where I performed system analysis, documentation, and user support.\",\n    \"I hold a B.Tech degree from NIT Durgapur and I am currently pursuing a Full Stack Data Science Bootcamp at iNeuron. I have completed coursework in machine learning, probability and statistics, calculus, and stock market analysis. I am proficient in Python, SQL, NLP, time series forecasting, deep learning, and various data science libraries such as NumPy, Pandas, TensorFlow, and Keras.\",\n    \"In addition to my technical skills, I possess strong soft skills including time management, teamwork, problem-solving, and presentation skills. I am a hardworking and flexible individual who is always eager to learn and adapt to new challenges.\"\n)\n\nforeach ($old in $oldParasToRemove) {\n    $r = Find-ExactRange $d $old\n    [void]$r.MoveEnd(1, 2)  # wdCharacter = 1; extend by 2 chars to swallow the trailing blank-line break pair\n    $r.Text = \"\"\n}\n"}
